# A new "September" notification entry ("exclusive on axis") was logged
# at 2024-09-24 17:28:35. It lands at the top of the existing September
# details list (row 48 of the "2024" sheet), pushing every row below it
# down by one (old row 48 -> new row 49, ..., old row 205 -> new row 206).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# Insert a new blank row at row 48, shifting rows 48:205 down to 49:206.
$ws.Rows("48:48").Insert()

# Populate the new row's September_Details / September_Date cells.
$ws.Range("R48").Value = "exclusive on axis"
$ws.Range("S48").Value = "2024-09-24 17:28:35"
